$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()

# Note text rewording
$wsAbout.Range("A10").Value = "Fuel use by fuel (fractions) are provided for CHP systems."

# Turn the IEA URL cell into a real hyperlink (Excel auto-applies the
# built-in "Hyperlink" style: underline, theme color 4, sz 10)
$wsAbout.Hyperlinks.Add($wsAbout.Range("B6"), $wsAbout.Range("B6").Value2) | Out-Null

$wsAbout.Range("A11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Data"
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("A1:H14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "BFoHPbF"
# ---------------------------------------------------------------------------
$wsBFoHPbF = $wb.Worksheets.Item("BFoHPbF")
$wsBFoHPbF.Activate()

# Widen column A and give the new header a title + wrap + bold + taller row
$wsBFoHPbF.Columns.Item(1).ColumnWidth = 23.3072916667
$wsBFoHPbF.Range("A1").Value = "Fraction of Heat Provided by Fuel (dimensionless)"
$wsBFoHPbF.Range("A1").Font.Bold = $true
$wsBFoHPbF.Range("A1").WrapText = $true
$wsBFoHPbF.Rows.Item(1).RowHeight = 30

# New fuel rows: crude oil / heavy-or-residual fuel oil / LPG / hydrogen
# Each mirrors the existing rows 2-7 pattern: B has a literal 0, and the
# rest of the row (C:AK) repeats "=$B<row>" across all year columns.
$wsBFoHPbF.Range("A8").Value = "crude oil"
$wsBFoHPbF.Range("B8").Value = 0
$wsBFoHPbF.Range("C8:AK8").FormulaR1C1 = "=RC2"

$wsBFoHPbF.Range("A9").Value = "heavy or residual fuel oil"
$wsBFoHPbF.Range("B9").Value = 0
$wsBFoHPbF.Range("C9:AK9").FormulaR1C1 = "=RC2"

$wsBFoHPbF.Range("A10").Value = "LPG propane or butane"
$wsBFoHPbF.Range("B10").Value = 0
$wsBFoHPbF.Range("C10:AK10").FormulaR1C1 = "=RC2"

$wsBFoHPbF.Range("A11").Value = "hydrogen"
$wsBFoHPbF.Range("B11").Value = 0
$wsBFoHPbF.Range("C11:AK11").FormulaR1C1 = "=RC2"

# Restore the originally-active sheet so tabSelected stays on "About"
$wsAbout.Activate()
